$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to remain as plain text so that
# numeric-looking values (e.g. "1.000", "0.7758") keep their original
# formatting instead of being coerced into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.794.51"
$ws.Range("E2").Value = "  -1.62%  "

# Row 3
$ws.Range("D3").Value = "1.891.90"
$ws.Range("E3").Value = "  -1.29%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "0.7758"
$ws.Range("E5").Value = "  -4.87%  "

# Row 6
$ws.Range("D6").Value = "243.76"
$ws.Range("E6").Value = "  -0.43%  "

# Row 7
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").Value = "0.3119"
$ws.Range("E8").Value = "  -4.16%  "

# Row 9
$ws.Range("D9").Value = "25.24"
$ws.Range("E9").Value = "  -7.22%  "

# Row 10
$ws.Range("D10").Value = "0.07157"
$ws.Range("E10").Value = "  +0.68%  "

# Row 11
$ws.Range("D11").Value = "0.08070"
$ws.Range("E11").Value = "  -0.23%  "

# Row 12
$ws.Range("D12").Value = "0.7644"
$ws.Range("E12").Value = "  -2.18%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.453"
$ws.Range("E13").Value = "  +1.50%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.871.78"
$ws.Range("E14").Value = "  -2.04%  "

# Row 15
$ws.Range("D15").Value = "92.23"
$ws.Range("E15").Value = "  -2.62%  "

# Row 16
$ws.Range("D16").Value = "6.153"
$ws.Range("E16").Value = "  +2.04%  "

# Row 17
$ws.Range("D17").Value = "29.821.16"
$ws.Range("E17").Value = "  -1.56%  "

# Row 18
$ws.Range("D18").Value = "13.94"
$ws.Range("E18").Value = "  -2.95%  "

# Row 19
$ws.Range("D19").Value = "243.47"
$ws.Range("E19").Value = "  -2.08%  "

# Row 20
$ws.Range("D20").Value = "0.000007752"
$ws.Range("E20").Value = "  -1.09%  "

# Row 21
$ws.Range("D21").Value = "1.000"

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.138.18"
$ws.Range("E22").Value = "  -1.93%  "

# Row 23
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "8.098"
$ws.Range("E23").Value = "  +5.61%  "

# Row 24
$ws.Range("E24").Value = "  -0.03%  "

# Row 25
$ws.Range("D25").Value = "0.1594"
$ws.Range("E25").Value = "  -3.32%  "

# Row 26
$ws.Range("D26").Value = "9.375"
$ws.Range("E26").Value = "  -0.81%  "

# Row 27
$ws.Range("D27").Value = "161.85"
$ws.Range("E27").Value = "  -3.71%  "

# Row 28
$ws.Range("D28").Value = "18.73"
$ws.Range("E28").Value = "  -2.14%  "

# Row 29
$ws.Range("D29").Value = "2.049"
$ws.Range("E29").Value = "  -3.60%  "

# Row 30
$ws.Range("D30").Value = "1.440"
$ws.Range("E30").Value = "  +4.97%  "

# Row 31
$ws.Range("E31").Value = "  +0.69%  "

# Row 32
$ws.Range("D32").Value = "4.458"
$ws.Range("E32").Value = "  +2.69%  "

# Row 33
$ws.Range("D33").Value = "4.095"
$ws.Range("E33").Value = "  -0.89%  "

# Row 34
$ws.Range("D34").Value = "0.05522"
$ws.Range("E34").Value = "  -3.03%  "

# Row 35
$ws.Range("D35").Value = "1.260"
$ws.Range("E35").Value = "  -3.27%  "

# Row 36
$ws.Range("D36").Value = "0.7447"
$ws.Range("E36").Value = "  +0.42%  "

# Row 37
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("D38").Value = "2.620"
$ws.Range("E38").Value = "  -3.67%  "

# Row 39
$ws.Range("D39").Value = "0.01912"
$ws.Range("E39").Value = "  -1.88%  "

# Row 40
$ws.Range("D40").Value = "2.777"
$ws.Range("E40").Value = "  -1.58%  "

# Row 41
$ws.Range("D41").Value = "1.135.32"
$ws.Range("E41").Value = "  +8.30%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.4418"
$ws.Range("E42").Value = "  -1.43%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "73.46"
$ws.Range("E43").Value = "  -0.78%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.867"
$ws.Range("E44").Value = "  -2.25%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "0.8510"
$ws.Range("E45").Value = "  -0.27%  "

# Row 46
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "103.96"
$ws.Range("E46").Value = "  +1.05%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  -0.12%  "

# Row 48
$ws.Range("D48").Value = "1.885"
$ws.Range("E48").Value = "  -2.44%  "

# Row 49
$ws.Range("D49").Value = "9.893"
$ws.Range("E49").Value = "  -0.73%  "

# Row 50
$ws.Range("D50").Value = "7.432"
$ws.Range("E50").Value = "  -2.48%  "

# Row 51
$ws.Range("D51").Value = "3.015"
$ws.Range("E51").Value = "  +10.19%  "
